$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 156791.03
$ws.Range("I132").Value = 162462.69
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 487388.07
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -484858.07
$ws.Range("N132").Value = -50060

# Row 135
$ws.Range("H135").Value = 800
$ws.Range("I135").Value = 800
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7200
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4665
$ws.Range("N135").ClearContents()

# Row 137
$ws.Range("H137").Value = 3623.5908
$ws.Range("I137").Value = 3088.2222
$ws.Range("J137").Value = 6032.75
$ws.Range("K137").Value = 9264.6666
$ws.Range("L137").Value = 18098.25
$ws.Range("M137").Value = -6714.6666
$ws.Range("N137").Value = -23198.25

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# Row 5
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 60
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = 52
$ws.Range("N5").Value = -284

# Row 61
$ws.Range("H61").Value = 2087.92
$ws.Range("I61").Value = 1619.1765
$ws.Range("J61").Value = 3084
$ws.Range("K61").Value = 1619.1765
$ws.Range("L61").Value = 3084
$ws.Range("M61").Value = -1407.1765
$ws.Range("N61").Value = -3508

# Row 122
$ws.Range("H122").Value = 2900.5625
$ws.Range("I122").Value = 1543.8
$ws.Range("J122").Value = 3517.2727
$ws.Range("K122").Value = 4631.4
$ws.Range("L122").Value = 10551.8181
$ws.Range("M122").Value = -2181.4
$ws.Range("N122").Value = -15451.8181

# Row 136
$ws.Range("H136").Value = 2087.92
$ws.Range("I136").Value = 1619.1765
$ws.Range("J136").Value = 3084
$ws.Range("K136").Value = 4857.529500000001
$ws.Range("L136").Value = 9252
$ws.Range("M136").Value = -2307.529500000001
$ws.Range("N136").Value = -14352

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 60
$ws.Range("L4").Value = 60
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = -290

# Row 134
$ws.Range("H134").Value = 3082.5918
$ws.Range("I134").Value = 1831.091
$ws.Range("J134").Value = 4102.3335
$ws.Range("K134").Value = 5493.272999999999
$ws.Range("L134").Value = 12307.0005
$ws.Range("M134").Value = -2958.272999999999
$ws.Range("N134").Value = -17377.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 22731328
$ws.Range("I31").Value = 2161.3076
$ws.Range("K31").Value = 2161.3076
$ws.Range("M31").Value = -1866.3076

# Row 34
$ws.Range("H34").Value = 22731328
$ws.Range("I34").Value = 2161.3076
$ws.Range("K34").Value = 2161.3076
$ws.Range("M34").Value = -1959.3076

# Row 82
$ws.Range("H82").Value = 39600
$ws.Range("J82").Value = 39600
$ws.Range("L82").Value = 39600
$ws.Range("N82").Value = -40322

# Row 85
$ws.Range("H85").Value = 39600
$ws.Range("J85").Value = 39600
$ws.Range("L85").Value = 39600
$ws.Range("N85").Value = -42096

# Row 134
$ws.Range("H134").Value = 7489.952
$ws.Range("I134").Value = 11151.9
$ws.Range("J134").Value = 4160.909
$ws.Range("K134").Value = 33455.7
$ws.Range("L134").Value = 12482.727
$ws.Range("M134").Value = -30920.7
$ws.Range("N134").Value = -17552.727

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 5501.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5501.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 16504.5
$ws.Range("N86").Value = -18876.5
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 5501.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5501.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 49513.5
$ws.Range("N89").Value = -61369.5
$ws.Range("M89").ClearContents()

# Row 113
$ws.Range("H113").Value = 561.0357
$ws.Range("I113").Value = 530.1111
$ws.Range("K113").Value = 1590.3333
$ws.Range("M113").Value = 579.6667000000002

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 1008000
$ws.Range("I3").Value = 1669333.6
$ws.Range("J3").Value = 15999.5
$ws.Range("K3").Value = 1669333.6
$ws.Range("L3").Value = 15999.5
$ws.Range("M3").Value = -1669217.6
$ws.Range("N3").Value = -16231.5

# Row 11
$ws.Range("H11").Value = 6779978
$ws.Range("I11").Value = 9900030
$ws.Range("J11").Value = 2879913
$ws.Range("K11").Value = 9900030
$ws.Range("L11").Value = 2879913
$ws.Range("M11").Value = -9899891
$ws.Range("N11").Value = -2880191

# Row 122
$ws.Range("H122").Value = 4743.4546
$ws.Range("I122").Value = 2771.25
$ws.Range("J122").Value = 10002.667
$ws.Range("K122").Value = 8313.75
$ws.Range("L122").Value = 30008.001
$ws.Range("M122").Value = -5863.75
$ws.Range("N122").Value = -34908.001

$ws = $wb.Worksheets.Item("LTW")
# Row 81
$ws.Range("H81").Value = 66324.5
$ws.Range("J81").Value = 66324.5
$ws.Range("L81").Value = 66324.5
$ws.Range("N81").Value = -68320.5

# Row 84
$ws.Range("H84").Value = 66324.5
$ws.Range("J84").Value = 66324.5
$ws.Range("L84").Value = 198973.5
$ws.Range("N84").Value = -208957.5

# Row 122
$ws.Range("H122").Value = 5086.706
$ws.Range("I122").Value = 2816.2727
$ws.Range("J122").Value = 9249.166999999999
$ws.Range("K122").Value = 8448.8181
$ws.Range("L122").Value = 27747.501
$ws.Range("M122").Value = -5998.8181
$ws.Range("N122").Value = -32647.501

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5495412
$ws.Range("I81").Value = 6494323.5
$ws.Range("J81").Value = 1400
$ws.Range("K81").Value = 12988647
$ws.Range("L81").Value = 2800
$ws.Range("M81").Value = -12987586
$ws.Range("N81").Value = -4922

# Row 84
$ws.Range("H84").Value = 5495412
$ws.Range("I84").Value = 6494323.5
$ws.Range("J84").Value = 1400
$ws.Range("K84").Value = 64943235
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = -64937931
$ws.Range("N84").Value = -24608
